$d = $word.ActiveDocument

# --- Simple, unambiguous whole-document replacements (each old string is
# either unique, or all of its occurrences should receive the same
# replacement text) ---

$d.Content.Find.Execute("英文", $true, $false, $false, $false, $false, $true, 1, $false, "英语", 2) | Out-Null

$d.Content.Find.Execute("葡萄牙文 / 法文 / 泰文 / 越南文 / 西班牙文", $true, $false, $false, $false, $false, $true, 1, $false, "葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语", 2) | Out-Null

$d.Content.Find.Execute("簡介", $true, $false, $false, $false, $false, $true, 1, $false, "简要", 2) | Out-Null

$d.Content.Find.Execute("發送給確認參加活動的參與者的電子郵件。 我們想與他們分享航班和住宿的預定詳情。", $true, $false, $false, $false, $false, $true, 1, $false, "发送给确认参加活动的与会者的电子邮件。 我们希望与他们分享航班和住宿的预订细节。", 2) | Out-Null

$d.Content.Find.Execute("目標受眾", $true, $false, $false, $false, $false, $true, 1, $false, "目标受众", 2) | Out-Null

$d.Content.Find.Execute("活動參加者", $true, $false, $false, $false, $false, $true, 1, $false, "活动与会者", 2) | Out-Null

$d.Content.Find.Execute("主題: ", $true, $false, $false, $false, $false, $true, 1, $false, "主题: ", 2) | Out-Null

$d.Content.Find.Execute("這裡是您的", $true, $false, $false, $false, $false, $true, 1, $false, "这是您关于", 2) | Out-Null

$d.Content.Find.Execute("[活動名稱]", $true, $false, $false, $false, $false, $true, 1, $false, "[活动名称]", 2) | Out-Null

$d.Content.Find.Execute(" 預訂詳情 ", $true, $false, $false, $false, $false, $true, 1, $false, "的预订详情 ", 2) | Out-Null

$d.Content.Find.Execute("我們迫不及待想見到您！ ", $true, $false, $false, $false, $false, $true, 1, $false, "我们迫不及待想与您会面！ ", 2) | Out-Null

$d.Content.Find.Execute("[合作夥伴姓名]", $true, $false, $false, $false, $false, $true, 1, $false, "[合作伙伴姓名]", 2) | Out-Null

$d.Content.Find.Execute("， ", $true, $false, $false, $false, $false, $true, 1, $false, ", ", 2) | Out-Null

$d.Content.Find.Execute("希望您和我們一樣對 ", $true, $false, $false, $false, $false, $true, 1, $false, "希望您和我们一样对 ", 2) | Out-Null

$d.Content.Find.Execute(" 感到興奮。 活動即將開始，我們已經做好了一切準備，希望您能和我們一起參加這次 ", $true, $false, $false, $false, $false, $true, 1, $false, " 感到兴奋。 活动即将开始，我们已经做好了一切准备，希望您能和我们一起参加这次 ", 2) | Out-Null

$d.Content.Find.Execute("會議/研討會/旅行", $true, $false, $false, $false, $false, $true, 1, $false, "会议/研讨会/旅行", 2) | Out-Null

$d.Content.Find.Execute("在此電子郵件中，我們已連結/附上以下文件：", $true, $false, $false, $false, $false, $true, 1, $false, "在这封电子邮件中，我们链接/附上了以下文件：", 2) | Out-Null

$d.Content.Find.Execute("您的來回機票", $true, $false, $false, $false, $false, $true, 1, $false, "您的来回机票", 2) | Out-Null

$d.Content.Find.Execute("您的住宿預訂詳情", $true, $false, $false, $false, $false, $true, 1, $false, "您的住宿预订详情", 2) | Out-Null

$d.Content.Find.Execute("您的簽證資料 ", $true, $false, $false, $false, $false, $true, 1, $false, "您的签证信息 ", 2) | Out-Null

$d.Content.Find.Execute("（如適用）", $true, $false, $false, $false, $false, $true, 1, $false, "（如适用）", 2) | Out-Null

$d.Content.Find.Execute("如有任何疑問，請與我們聯繫：", $true, $false, $false, $false, $false, $true, 1, $false, "如果有任何疑问，请联系我们：", 2) | Out-Null

$d.Content.Find.Execute("即時聊天", $true, $false, $false, $false, $false, $true, 1, $false, "实时聊天", 2) | Out-Null

# --- Targeted, position-scoped replacements. Several of the remaining old
# strings (" 或 ", "WhatsApp", "。 ", "在") occur more than once in the body,
# but only specific occurrences change here, so each Find is chained:
# it starts searching right where the previous replacement ended, walking
# forward through the document in order. ---

$r = $d.Content
$r.Start = 0
$r.Find.Execute(" 或 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = " |"
$pos = $r.End

$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute("WhatsApp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = " WhatsApp"
$pos = $r.End

$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute("。 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = ". "
$pos = $r.End

$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute("如有任何疑問，請通過 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = "如有任何疑问，请通过 "
$pos = $r.End

$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute("[電子郵件地址]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = "[电子邮件地址]"
$pos = $r.End

# the second " 或 " (between e-mail and WHATSAPP number placeholders) is
# left unchanged, so just skip past it to keep the chain positioned
$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute(" 或 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos = $r.End

$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute("[WHATSAPP 號碼]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = "[WHATSAPP 号码]"
$pos = $r.End

$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute(" (WhatsApp) 聯繫您的區域經理 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = " (WhatsApp) 联系您的区域经理 "
$pos = $r.End

$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute("[NAME]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = "[姓名]"
$pos = $r.End

$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute("在", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = "在 "
$pos = $r.End

$r = $d.Range($pos, $d.Content.End)
$r.Find.Execute("見！", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = " 见！"

# --- Comments: replace the whole text of each comment. Comments.Item() is
# ordered by position of the comment reference in the body (1 = id 0,
# 2 = id 1, 3 = id 2). Direct Range.Text assignment is used because
# Find does not resolve comment story ranges in this object model. ---

$d.Comments.Item(1).Range.Text = "选择一个"
$d.Comments.Item(2).Range.Text = "检查这些是否是包含的文件"
$d.Comments.Item(3).Range.Text = "选择其一"
